$wb = $excel.ActiveWorkbook

$timestamp = "2025-12-19 03:03:23"

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $timestamp
    }
}
